$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) keeps its text format, since some values
# (e.g. "0.998", "211.19") would otherwise be auto-converted to numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '27.869.32'
$ws.Range("E2").Value = '  -0.25%  '

$ws.Range("D3").Value = '1.628.26'
$ws.Range("E3").Value = '  -0.78%  '

$ws.Range("D4").Value = '0.998'
$ws.Range("E4").Value = '  -0.20%  '

$ws.Range("D5").Value = '211.19'
$ws.Range("E5").Value = '  -0.70%  '

$ws.Range("D6").Value = '0.522'
$ws.Range("E6").Value = '  -0.47%  '

$ws.Range("E7").Value = '  -0.21%  '

$ws.Range("D8").Value = '23.32'
$ws.Range("E8").Value = '  -0.49%  '

$ws.Range("D9").Value = '0.257'
$ws.Range("E9").Value = '  -0.56%  '

$ws.Range("D10").Value = '0.0612'
$ws.Range("E10").Value = '  -0.42%  '

$ws.Range("D11").Value = '0.0880'
$ws.Range("E11").Value = '  -0.39%  '

$ws.Range("D12").Value = '1.855.11'
$ws.Range("E12").Value = '  -1.00%  '

$ws.Range("D13").Value = '1.617.09'
$ws.Range("E13").Value = '  -1.46%  '

$ws.Range("E14").Value = '  -1.47%  '

$ws.Range("D15").Value = '0.562'
$ws.Range("E15").Value = '  -1.51%  '

$ws.Range("D16").Value = '65.23'
$ws.Range("E16").Value = '  -0.37%  '

$ws.Range("D17").Value = '27.856.09'
$ws.Range("E17").Value = '  -0.27%  '

$ws.Range("D18").Value = '229.23'
$ws.Range("E18").Value = '  -1.59%  '

$ws.Range("D19").Value = '7.67'
$ws.Range("E19").Value = '  +0.91%  '

$ws.Range("D20").Value = '0.0₃0720'
$ws.Range("E20").Value = '  -0.30%  '

$ws.Range("D21").Value = '0.997'
$ws.Range("E21").Value = '  -0.31%  '

$ws.Range("D22").Value = '4.33'
$ws.Range("E22").Value = '  -1.02%  '

$ws.Range("D23").Value = '10.08'
$ws.Range("E23").Value = '  -3.38%  '

$ws.Range("E24").Value = '  -2.28%  '

$ws.Range("D25").Value = '154.09'
$ws.Range("E25").Value = '  +0.65%  '

$ws.Range("D26").Value = '6.89'
$ws.Range("E26").Value = '  -0.01%  '

$ws.Range("E27").Value = '  -0.27%  '

$ws.Range("D28").Value = '15.52'
$ws.Range("E28").Value = '  -0.92%  '

$ws.Range("D29").Value = '0.998'
$ws.Range("E29").Value = '  -0.24%  '

$ws.Range("D30").Value = '1.18'
$ws.Range("E30").Value = '  -0.96%  '

$ws.Range("D31").Value = '0.0481'
$ws.Range("E31").Value = '  -0.71%  '

$ws.Range("D32").Value = '3.41'
$ws.Range("E32").Value = '  +0.25%  '

$ws.Range("D33").Value = '3.10'
$ws.Range("E33").Value = '  -0.13%  '

$ws.Range("D34").Value = '1.393.00'
$ws.Range("E34").Value = '  -1.03%  '

$ws.Range("D35").Value = '1.59'
$ws.Range("E35").Value = '  +0.25%  '

$ws.Range("E36").Value = '  +10.00%  '

$ws.Range("E37").Value = '  -1.09%  '

$ws.Range("E38").Value = '  +0.20%  '

$ws.Range("D39").Value = '0.557'
$ws.Range("E39").Value = '  -1.09%  '

$ws.Range("D40").Value = '0.851'
$ws.Range("E40").Value = '  -3.21%  '

$ws.Range("D41").Value = '0.997'
$ws.Range("E41").Value = '  -0.27%  '

$ws.Range("E42").Value = '  -1.10%  '

$ws.Range("E43").Value = '  -0.18%  '

$ws.Range("D44").Value = '65.73'
$ws.Range("E44").Value = '  -2.19%  '

$ws.Range("D45").Value = '5.42'
$ws.Range("E45").Value = '  -1.68%  '

$ws.Range("D46").Value = '1.765.13'
$ws.Range("E46").Value = '  -1.03%  '

$ws.Range("E47").Value = '  -2.81%  '

$ws.Range("D48").Value = '88.03'
$ws.Range("E48").Value = '  +0.27%  '

$ws.Range("D49").Value = '0.102'
$ws.Range("E49").Value = '  +1.60%  '

$ws.Range("D50").Value = '0.0503'
$ws.Range("E50").Value = '  -0.57%  '

$ws.Range("D51").Value = '7.60'
$ws.Range("E51").Value = '  +0.24%  '
